$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("test" task) is fleshed out into a real Herbie usability test entry
$ws.Range("B4").Value = "herbie_basic_test"
$ws.Range("C4").Value = "enter firstname last name and click on submit verify whther its correct or not"

# D4 gets a url value that is then turned into a clickable hyperlink
$ws.Range("D4").Value = "https://mieweb.github.io/herbie/playgrounds/login.html"
$ws.Hyperlinks.Add($ws.Range("D4"), "https://mieweb.github.io/herbie/playgrounds/login.html") | Out-Null

# E4 holds the herbie verification script, wrapped like the other script cells
$ws.Range("E4").Value = 'verify state is visible in "result"'
$ws.Range("E4").WrapText = $true

# Selection ends up on the newly completed cell
$ws.Range("E4").Select()
